$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant values (headless COM host doesn't predefine the xl* enums)
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlDouble = -4119
$xlThin = 1

# ---------------------------------------------------------------------
# 1) Shift the existing A:D data to B:E by inserting a new column at A
#    (the new first column becomes the "Task Id" column, matching the
#    other columns already carrying Task / Due Date / Status).
# ---------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# Column widths for the two brand-new columns (B = Task Id, E = Status).
$ws.Columns("B").ColumnWidth = 7.333333333333334
$ws.Columns("E").ColumnWidth = 12.666666666666668

# ---------------------------------------------------------------------
# 2) New Status values
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "Done"
$ws.Range("E4").Value = "Done"
$ws.Range("E5").Value = "Not Required"

# ---------------------------------------------------------------------
# 3) Header row (row 2) styling: bold white text on a black fill, with a
#    double rule around the outside of the header band.
# ---------------------------------------------------------------------
$header = $ws.Range("B2:E2")
$header.Font.Bold = $true
$header.Font.ThemeColor = 2          # xlThemeColorLight1 -> theme="0" (white)
$header.Interior.ThemeColor = 1      # xlThemeColorDark1  -> theme="1" (black)

$header.Borders.Item($xlEdgeTop).LineStyle = $xlDouble
$header.Borders.Item($xlEdgeBottom).LineStyle = $xlThin
$ws.Range("B2").Borders.Item($xlEdgeLeft).LineStyle = $xlDouble
$ws.Range("E2").Borders.Item($xlEdgeRight).LineStyle = $xlDouble
$ws.Range("C2:E2").Borders.Item($xlEdgeLeft).LineStyle = $xlThin
$ws.Range("B2:D2").Borders.Item($xlEdgeRight).LineStyle = $xlThin

# ---------------------------------------------------------------------
# 4) Body rows (3-16): thin grid inside the table, double rule on the
#    outer left/right edges of the whole block.
# ---------------------------------------------------------------------
$body = $ws.Range("B3:E16")
$body.Borders.Item($xlEdgeTop).LineStyle = $xlThin
$body.Borders.Item($xlEdgeBottom).LineStyle = $xlThin
$ws.Range("B3:B16").Borders.Item($xlEdgeLeft).LineStyle = $xlDouble
$ws.Range("C3:E16").Borders.Item($xlEdgeLeft).LineStyle = $xlThin
$ws.Range("B3:D16").Borders.Item($xlEdgeRight).LineStyle = $xlThin
$ws.Range("E3:E16").Borders.Item($xlEdgeRight).LineStyle = $xlDouble

# Status fills: green "Done", tinted-orange "Not Required".
$ws.Range("E3:E4").Interior.Color = 5296274      # BGR packing of FF92D050
$ws.Range("E5").Interior.ThemeColor = 10         # xlThemeColorAccent6 -> theme="9"

# ---------------------------------------------------------------------
# 5) Last data row (17): thin grid on top, double rule closing the
#    bottom of the whole block.
# ---------------------------------------------------------------------
$last = $ws.Range("B17:E17")
$last.Borders.Item($xlEdgeTop).LineStyle = $xlThin
$last.Borders.Item($xlEdgeBottom).LineStyle = $xlDouble
$ws.Range("B17").Borders.Item($xlEdgeLeft).LineStyle = $xlDouble
$ws.Range("C17:E17").Borders.Item($xlEdgeLeft).LineStyle = $xlThin
$ws.Range("B17:D17").Borders.Item($xlEdgeRight).LineStyle = $xlThin
$ws.Range("E17").Borders.Item($xlEdgeRight).LineStyle = $xlDouble

# ---------------------------------------------------------------------
# 6) Frame rows above/below the table (new rows 1 and 18) so the double
#    border reads as a closed box.
# ---------------------------------------------------------------------
$ws.Rows("1").RowHeight = 15.75
$ws.Rows("18").RowHeight = 15.75

# ---------------------------------------------------------------------
# 7) Misc sheet-level bits
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait
$ws.Range("E6").Select()

Write-Output "layout complete"
